$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 149.5
$ws.Range("I2").Value = 99
$ws.Range("J2").Value = 200
$ws.Range("K2").Value = 99
$ws.Range("L2").Value = 200
$ws.Range("M2").Value = 14
$ws.Range("N2").Value = -426
$ws.Range("H9").Value = 60.727272
$ws.Range("I9").Value = 60.875
$ws.Range("K9").Value = 60.875
$ws.Range("M9").Value = 108.125
$ws.Range("H32").Value = 2068.923
$ws.Range("I32").Value = 2211.111
$ws.Range("J32").Value = 1749
$ws.Range("K32").Value = 2211.111
$ws.Range("L32").Value = 1749
$ws.Range("M32").Value = -1885.111
$ws.Range("N32").Value = -2401
$ws.Range("H43").Value = 5993
$ws.Range("I43").Value = 5993
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 5993
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -5924
$ws.Range("N43").ClearContents()
$ws.Range("H58").Value = 1627.6666
$ws.Range("I58").Value = 787.25
$ws.Range("J58").Value = 3308.5
$ws.Range("K58").Value = 2361.75
$ws.Range("L58").Value = 9925.5
$ws.Range("M58").Value = -2211.75
$ws.Range("N58").Value = -10225.5
$ws.Range("H86").Value = 17250.25
$ws.Range("I86").Value = 17001
$ws.Range("J86").Value = 17499.5
$ws.Range("K86").Value = 17001
$ws.Range("L86").Value = 17499.5
$ws.Range("M86").Value = -15878
$ws.Range("N86").Value = -19745.5
$ws.Range("H88").Value = 3828.353
$ws.Range("I88").Value = 5187.25
$ws.Range("J88").Value = 2620.4443
$ws.Range("K88").Value = 5187.25
$ws.Range("L88").Value = 2620.4443
$ws.Range("M88").Value = -4781.25
$ws.Range("N88").Value = -3432.4443
$ws.Range("H89").Value = 17250.25
$ws.Range("I89").Value = 17001
$ws.Range("J89").Value = 17499.5
$ws.Range("K89").Value = 85005
$ws.Range("L89").Value = 87497.5
$ws.Range("M89").Value = -79389
$ws.Range("N89").Value = -98729.5
$ws.Range("H91").Value = 3828.353
$ws.Range("I91").Value = 5187.25
$ws.Range("J91").Value = 2620.4443
$ws.Range("K91").Value = 5187.25
$ws.Range("L91").Value = 2620.4443
$ws.Range("M91").Value = -3783.25
$ws.Range("N91").Value = -5428.4443
$ws.Range("H132").Value = 21083.7
$ws.Range("I132").Value = 21083.7
$ws.Range("K132").Value = 63251.10000000001
$ws.Range("M132").Value = -60721.10000000001
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2264.6296
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 2264.6296
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 2264.6296
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -3076.6296
$ws.Range("H91").Value = 2264.6296
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 2264.6296
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 2264.6296
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -5072.6296
$ws.Range("H110").Value = 2427.7856
$ws.Range("I110").Value = 3497.25
$ws.Range("K110").Value = 3497.25
$ws.Range("M110").Value = -1452.25
$ws.Range("H132").Value = 5363.0386
$ws.Range("I132").Value = 4247.448
$ws.Range("J132").Value = 6769.6523
$ws.Range("K132").Value = 12742.344
$ws.Range("L132").Value = 20308.9569
$ws.Range("M132").Value = -10212.344
$ws.Range("N132").Value = -25368.9569

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 6373.5713
$ws.Range("I105").Value = 935.6
$ws.Range("K105").Value = 935.6
$ws.Range("M105").Value = 811.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 96.28
$ws.Range("I7").Value = 65.21429000000001
$ws.Range("K7").Value = 65.21429000000001
$ws.Range("M7").Value = 47.78570999999999
$ws.Range("H16").Value = 141406.25
$ws.Range("I16").Value = 25650
$ws.Range("J16").Value = 334333.34
$ws.Range("K16").Value = 25650
$ws.Range("L16").Value = 334333.34
$ws.Range("M16").Value = -25363
$ws.Range("N16").Value = -334907.34
$ws.Range("H31").Value = 2097.889
$ws.Range("I31").Value = 2442.1765
$ws.Range("K31").Value = 2442.1765
$ws.Range("M31").Value = -2147.1765
$ws.Range("H34").Value = 2097.889
$ws.Range("I34").Value = 2442.1765
$ws.Range("K34").Value = 2442.1765
$ws.Range("M34").Value = -2240.1765
$ws.Range("H113").Value = 141406.25
$ws.Range("I113").Value = 25650
$ws.Range("J113").Value = 334333.34
$ws.Range("K113").Value = 25650
$ws.Range("L113").Value = 334333.34
$ws.Range("M113").Value = -23480
$ws.Range("N113").Value = -338673.34
$ws.Range("H132").Value = 5741.815
$ws.Range("I132").Value = 6214.9546
$ws.Range("J132").Value = 3660
$ws.Range("K132").Value = 18644.8638
$ws.Range("L132").Value = 10980
$ws.Range("M132").Value = -16114.8638
$ws.Range("N132").Value = -16040
$ws.Range("H134").Value = 2736.0312
$ws.Range("I134").Value = 2629.8235
$ws.Range("J134").Value = 2856.4
$ws.Range("K134").Value = 7889.470499999999
$ws.Range("L134").Value = 8569.200000000001
$ws.Range("M134").Value = -5354.470499999999
$ws.Range("N134").Value = -13639.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 309.66666
$ws.Range("J75").Value = 309.66666
$ws.Range("L75").Value = 928.9999799999999
$ws.Range("N75").Value = -2924.99998
$ws.Range("H78").Value = 309.66666
$ws.Range("J78").Value = 309.66666
$ws.Range("L78").Value = 2786.99994
$ws.Range("N78").Value = -12770.99994
$ws.Range("H117").Value = 643.0909
$ws.Range("I117").Value = 182.57143
$ws.Range("K117").Value = 547.71429
$ws.Range("M117").Value = 2894.28571
$ws.Range("H138").Value = 23443.408
$ws.Range("I138").Value = 66911
$ws.Range("K138").Value = 200733
$ws.Range("M138").Value = -195593
$ws.Range("H139").Value = 4058
$ws.Range("I139").Value = 3028.3635
$ws.Range("K139").Value = 9085.0905
$ws.Range("M139").Value = -3945.0905
$ws.Range("H141").Value = 14500
$ws.Range("I141").Value = 8000
$ws.Range("K141").Value = 24000
$ws.Range("M141").Value = -18820

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1884.9166
$ws.Range("I16").Value = 1884.9166
$ws.Range("K16").Value = 1884.9166
$ws.Range("M16").Value = -1714.9166
$ws.Range("H55").Value = 1302.7222
$ws.Range("I55").Value = 1132.7333
$ws.Range("K55").Value = 1132.7333
$ws.Range("M55").Value = -959.7333000000001
$ws.Range("H61").Value = 9950.684999999999
$ws.Range("I61").Value = 8357.235000000001
$ws.Range("K61").Value = 8357.235000000001
$ws.Range("M61").Value = -8155.235000000001
$ws.Range("H113").Value = 9950.684999999999
$ws.Range("I113").Value = 8357.235000000001
$ws.Range("K113").Value = 8357.235000000001
$ws.Range("M113").Value = -6187.235000000001
$ws.Range("H132").Value = 3247.6667
$ws.Range("I132").Value = 2636.4443
$ws.Range("K132").Value = 7909.3329
$ws.Range("M132").Value = -5379.3329
$ws.Range("H136").Value = 10198.923
$ws.Range("J136").Value = 22556.75
$ws.Range("L136").Value = 67670.25
$ws.Range("N136").Value = -72770.25
